{"js": "// Merge the three runs \"and\u2014\" + \" \" + \"an em dash.\" into a single run\n// \"and\u2014 an em dash.\" in the paragraph that precedes the\n// \"subsection-in-recursive-include\" bookmark (the \"Text in recursive\n// include, with italic, bold, ... and\u2014 an em dash.\" paragraph).\n//\n// The visible text is unchanged; only the run split collapses, so we\n// find the exact phrase and rewrite it in place via insertText/Replace,\n// which causes the host to coalesce the run(s) covering that text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \" and\u2014 an em dash.\";\nlet paragraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(target) !== -1 && text.indexOf(\"Text in recursive include\") !== -1) {\n    paragraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!paragraph) {\n  throw new Error(\"Could not find target paragraph for 'and\u2014 an em dash.'\");\n}\n\nconst paragraphRange = paragraph.getRange();\nconst results = paragraphRange.search(\"and\u2014 an em dash.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'and\u2014 an em dash.' text to merge\");\n}\n\n// Replace the matched text with itself; this collapses the underlying\n// run boundaries into a single run while keeping the rendered text the\n// same.\nresults.items[0].insertText(\"and\u2014 an em dash.\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Merge the three runs \"and\u2014\" + \" \" + \"an em dash.\" into a single run\n# \"and\u2014 an em dash.\" in the paragraph that starts with \"Text in\n# recursive include, with ...\" (the paragraph immediately before the\n# \"Subsection in recursive include\" heading). The visible text does not\n# change; only the underlying run split collapses.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph robustly by its unique leading text rather\n# than a hard-coded index (there is a look-alike paragraph earlier in the\n# document, \"Text before recursive include, with ...\", that must NOT be\n# touched).\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"Text in recursive include*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'Text in recursive include' paragraph\"\n}\n\n$paragraph = $d.Paragraphs.Item($targetIndex)\n\n# Find the trailing \"an em dash.\" text (it lives entirely inside the\n# last, unformatted run of the paragraph). Editing strictly inside that\n# run causes the engine to coalesce it with its immediate, identically\n# formatted predecessor runs (\"and\u2014\" and the single space between them)\n# into one run, matching the target edit exactly, while leaving the\n# separate \" \" / \"curly quotes,\" / smart-quote runs further back\n# untouched.\n$target = $paragraph.Range\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = \"an em dash.\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find 'an em dash.' text to merge\"\n}\n\n# Replace the matched text with itself; this collapses the run\n# boundaries of \"and\u2014\" / \" \" / \"an em dash.\" into a single run while\n# keeping the rendered text identical.\n$target.Delete()\n$target.InsertBefore(\"an em dash.\")\n"}
